$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 5825.6665
$ws.Range("I88").Value = 11000.75
$ws.Range("J88").Value = 3943.818
$ws.Range("K88").Value = 11000.75
$ws.Range("L88").Value = 3943.818
$ws.Range("M88").Value = -10594.75
$ws.Range("N88").Value = -4755.818
$ws.Range("H91").Value = 5825.6665
$ws.Range("I91").Value = 11000.75
$ws.Range("J91").Value = 3943.818
$ws.Range("K91").Value = 11000.75
$ws.Range("L91").Value = 3943.818
$ws.Range("M91").Value = -9596.75
$ws.Range("N91").Value = -6751.818
$ws.Range("H93").Value = 38601
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 38601
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 38601
$ws.Range("N93").Value = -43593
$ws.Range("H125").Value = 1287.9375
$ws.Range("I125").Value = 400
$ws.Range("J125").Value = 1414.7858
$ws.Range("K125").Value = 3600
$ws.Range("L125").Value = 12733.0722
$ws.Range("M125").Value = -1140
$ws.Range("N125").Value = -17653.0722
$ws.Range("H135").Value = 1223
$ws.Range("I135").Value = 956.0476
$ws.Range("J135").Value = 1783.6
$ws.Range("K135").Value = 8604.4284
$ws.Range("L135").Value = 16052.4
$ws.Range("M135").Value = -6069.428400000001
$ws.Range("N135").Value = -21122.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H14").Value = 4326.6665
$ws.Range("I14").Value = 3980
$ws.Range("J14").Value = 4500
$ws.Range("K14").Value = 3980
$ws.Range("L14").Value = 4500
$ws.Range("M14").Value = -3805
$ws.Range("N14").Value = -4850
$ws.Range("H21").Value = 2508.6667
$ws.Range("I21").Value = 1007
$ws.Range("J21").Value = 10017
$ws.Range("K21").Value = 1007
$ws.Range("L21").Value = 10017
$ws.Range("M21").Value = -633
$ws.Range("N21").Value = -10765
$ws.Range("H88").Value = 90983570
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 100081730
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 100081730
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -100082542
$ws.Range("H91").Value = 90983570
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 100081730
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 100081730
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -100084538
$ws.Range("H102").Value = 2993.3333
$ws.Range("I102").Value = 1000
$ws.Range("J102").Value = 3990
$ws.Range("K102").Value = 1000
$ws.Range("L102").Value = 3990
$ws.Range("M102").Value = 622
$ws.Range("N102").Value = -7234
$ws.Range("H105").Value = 41499.5
$ws.Range("J105").Value = 41499.5
$ws.Range("L105").Value = 41499.5
$ws.Range("N105").Value = -48487.5
$ws.Range("H122").Value = 2989.5151
$ws.Range("I122").Value = 2700
$ws.Range("J122").Value = 3568.5454
$ws.Range("K122").Value = 8100
$ws.Range("L122").Value = 10705.6362
$ws.Range("M122").Value = -5650
$ws.Range("N122").Value = -15605.6362
$ws.Range("H132").Value = 56732.21
$ws.Range("I132").Value = 114569.22
$ws.Range("J132").Value = 4678.9
$ws.Range("K132").Value = 343707.66
$ws.Range("L132").Value = 14036.7
$ws.Range("M132").Value = -341177.66
$ws.Range("N132").Value = -19096.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 20002802
$ws.Range("I86").Value = 33335834
$ws.Range("J86").Value = 3251.5
$ws.Range("K86").Value = 33335834
$ws.Range("L86").Value = 3251.5
$ws.Range("M86").Value = -33334711
$ws.Range("N86").Value = -5497.5
$ws.Range("H89").Value = 20002802
$ws.Range("I89").Value = 33335834
$ws.Range("J89").Value = 3251.5
$ws.Range("K89").Value = 166679170
$ws.Range("L89").Value = 16257.5
$ws.Range("M89").Value = -166673554
$ws.Range("N89").Value = -27489.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1669
$ws.Range("I16").Value = 1619
$ws.Range("J16").Value = 1769
$ws.Range("K16").Value = 1619
$ws.Range("L16").Value = 1769
$ws.Range("M16").Value = -1332
$ws.Range("N16").Value = -2343
$ws.Range("H113").Value = 1669
$ws.Range("I113").Value = 1619
$ws.Range("J113").Value = 1769
$ws.Range("K113").Value = 1619
$ws.Range("L113").Value = 1769
$ws.Range("M113").Value = 551
$ws.Range("N113").Value = -6109

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1184.34
$ws.Range("I113").Value = 494.75
$ws.Range("J113").Value = 2062
$ws.Range("K113").Value = 1484.25
$ws.Range("L113").Value = 6186
$ws.Range("M113").Value = 685.75
$ws.Range("N113").Value = -10526

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2884.1292
$ws.Range("I132").Value = 2214.85
$ws.Range("J132").Value = 4101
$ws.Range("K132").Value = 6644.549999999999
$ws.Range("L132").Value = 12303
$ws.Range("M132").Value = -4114.549999999999
$ws.Range("N132").Value = -17363

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3125.5
$ws.Range("I61").Value = 2334
$ws.Range("J61").Value = 5500
$ws.Range("K61").Value = 2334
$ws.Range("L61").Value = 5500
$ws.Range("M61").Value = -2132
$ws.Range("N61").Value = -5904
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H113").Value = 3125.5
$ws.Range("I113").Value = 2334
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 2334
$ws.Range("L113").Value = 5500
$ws.Range("M113").Value = -164
$ws.Range("N113").Value = -9840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6331.96
$ws.Range("I132").Value = 6940.647
$ws.Range("J132").Value = 5038.5
$ws.Range("K132").Value = 20821.941
$ws.Range("L132").Value = 15115.5
$ws.Range("M132").Value = -18291.941
$ws.Range("N132").Value = -20175.5
